# Daily update at 8 AM UTC
# Appends the next day's win totals to the tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 45 was previously the last data row and carried the "date only"
# number format reserved for the final row. Since a new row now follows
# it, it reverts to the standard date+time format used by the other
# interior rows.
$ws.Range("A45").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data as the new last row, using the "date only"
# number format reserved for the final row.
$ws.Range("A46").Value = 45631
$ws.Range("A46").NumberFormat = "YYYY-MM-DD"
$ws.Range("B46").Value = 120
$ws.Range("C46").Value = 101
$ws.Range("D46").Value = 111
